# Reprocess rhyolitic glasses and update summary figure
# Update original.kraw_stdev_apf_pcnt (E), corrected.kraw_apf_pcnt (H),
# corrected.kraw_stdev_apf_pcnt (I), montecarlo.kraw_pcnt (J),
# montecarlo.kraw_stdev_pcnt (K), montecarlo.kraw_apf_pcnt (L) and
# montecarlo.kraw_stdev_apf_pcnt (M) for each spot row (2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - A876_50nA_20um_146.4_spot-1
$ws.Range("E2").Value = 9.720000000000001
$ws.Range("H2").Value = 4.58
$ws.Range("I2").Value = 9.720000000000001
$ws.Range("J2").Value = 3.8
$ws.Range("K2").Value = 2.88
$ws.Range("L2").Value = 4.58
$ws.Range("M2").Value = 4.62

# Row 3 - A876_50nA_20um_146.4_spot-2
$ws.Range("E3").Value = 10.47
$ws.Range("I3").Value = 10.47
$ws.Range("K3").Value = 3.57
$ws.Range("M3").Value = 5.08

# Row 4 - A876_50nA_20um_146.4_spot-3
$ws.Range("E4").Value = 13.19
$ws.Range("H4").Value = 4.18
$ws.Range("I4").Value = 13.19
$ws.Range("J4").Value = 3.47
$ws.Range("K4").Value = 4.02
$ws.Range("L4").Value = 4.18
$ws.Range("M4").Value = 5.41

# Row 5 - A876_50nA_20um_146.4_spot-4
$ws.Range("E5").Value = 9.76
$ws.Range("H5").Value = 4.77
$ws.Range("I5").Value = 9.76
$ws.Range("J5").Value = 3.96
$ws.Range("K5").Value = 3.42
$ws.Range("L5").Value = 4.77
$ws.Range("M5").Value = 4.98

# Row 6 - A876_50nA_20um_146.4_spot-5
$ws.Range("E6").Value = 11.78
$ws.Range("H6").Value = 4.14
$ws.Range("I6").Value = 11.78
$ws.Range("J6").Value = 3.44
$ws.Range("K6").Value = 3.18
$ws.Range("L6").Value = 4.14
$ws.Range("M6").Value = 4.81

# Row 7 - A876_50nA_20um_146.4_spot-6
$ws.Range("E7").Value = 236.78
$ws.Range("H7").Value = 2.81
$ws.Range("I7").Value = 236.78
$ws.Range("J7").Value = 2.33
$ws.Range("K7").Value = 4.95
$ws.Range("L7").Value = 2.81
$ws.Range("M7").Value = 6.13
